# diagnostic.xlsx edit: "disconnected_elements" contingency table
#  B1 = 0                  (bold, bordered, centered/top-aligned)
#  A2 = 0                  (same formatting as B1)
#  B2 = "disconnected_elements"   (plain, shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values ---------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- formatting for B1 (build the style once) ------------------------
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment   = -4160   # xlTop
$r1.Borders.LineStyle   = 1       # xlContinuous -> thin
$r1.Borders.ColorIndex  = -4105   # xlColorIndexAutomatic

# --- copy the exact same style onto A2 in one shot --------------------
# (applying the same properties again cell-by-cell would create stray,
#  unreferenced intermediate cellXfs/ border entries, so we clone the
#  already-built format instead)
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)           # xlPasteFormats
$r2.Value = 0
$excel.CutCopyMode = $false
